# Swap the full data (columns B:AD) between pairs of rows that were
# entered in the wrong order. Column A (the sequential row index) is
# left untouched on purpose.
#
# Pairs to swap:
#   rows 4  <-> 5
#   rows 164 <-> 165
#   rows 181 <-> 182
#   rows 190 <-> 191

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AD$rowA")
    $rangeB = $ws.Range("B$rowB`:AD$rowB")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

Swap-Rows 4 5
Swap-Rows 164 165
Swap-Rows 181 182
Swap-Rows 190 191
